$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Drop rows 113,114,116,118,123-131,133-134,136-138,140-146 entirely - clearing
# their content removes them from the sheetData once no cells remain.
$ws.Rows("113:114").ClearContents()
$ws.Rows("116").ClearContents()
$ws.Rows("118").ClearContents()
$ws.Rows("123:131").ClearContents()
$ws.Rows("133:134").ClearContents()
$ws.Rows("136:138").ClearContents()
$ws.Rows("140:146").ClearContents()

# Surviving rows keep their styled cells but lose their values.
$ws.Range("C115:L115").ClearContents()
$ws.Range("N115:W115").ClearContents()

$ws.Range("C117:L117").ClearContents()
$ws.Range("N117:W117").ClearContents()

$ws.Range("C119:L119").ClearContents()
$ws.Range("N119:W119").ClearContents()

$ws.Range("C120:L120").ClearContents()
$ws.Range("N120:W120").ClearContents()

$ws.Range("C121:L121").ClearContents()
$ws.Range("N121:W121").ClearContents()

$ws.Range("C122:L122").ClearContents()
$ws.Range("N122:W122").ClearContents()

$ws.Range("C132:L132").ClearContents()
$ws.Range("N132:W132").ClearContents()

$ws.Range("C135:L135").ClearContents()
$ws.Range("N135:W135").ClearContents()

$ws.Range("C139:L139").ClearContents()
$ws.Range("N139:W139").ClearContents()

# Move the active selection to M123, replacing the old N139:W139 selection.
$ws.Activate()
$ws.Range("M123").Select()
